$wb = $excel.ActiveWorkbook

# Map of sheet index -> { row -> new F value }, applying the "想去人数" (want-to-go count) updates

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1292
$ws.Range("F7").Value = 995
$ws.Range("F8").Value = 947
$ws.Range("F15").Value = 4228
$ws.Range("F16").Value = 1240
$ws.Range("F18").Value = 2710
$ws.Range("F20").Value = 1105
$ws.Range("F21").Value = 3704
$ws.Range("F22").Value = 796
$ws.Range("F24").Value = 47
$ws.Range("F25").Value = 1443
$ws.Range("F30").Value = 974
$ws.Range("F31").Value = 232
$ws.Range("F33").Value = 19
$ws.Range("F34").Value = 1404
$ws.Range("F35").Value = 1985
$ws.Range("F36").Value = 945
$ws.Range("F37").Value = 6
$ws.Range("F38").Value = 511
$ws.Range("F41").Value = 600
$ws.Range("F42").Value = 295
$ws.Range("F43").Value = 112
$ws.Range("F46").Value = 86

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 122

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 476

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 476
$ws.Range("F4").Value = 1292
$ws.Range("F6").Value = 995
$ws.Range("F7").Value = 947
$ws.Range("F16").Value = 4228
$ws.Range("F17").Value = 1240
$ws.Range("F20").Value = 2710
$ws.Range("F22").Value = 1105
$ws.Range("F23").Value = 3704
$ws.Range("F24").Value = 796
$ws.Range("F27").Value = 47
$ws.Range("F31").Value = 122
$ws.Range("F34").Value = 974
$ws.Range("F35").Value = 232
$ws.Range("F37").Value = 1404
$ws.Range("F38").Value = 1985
$ws.Range("F40").Value = 945
$ws.Range("F43").Value = 511
$ws.Range("F45").Value = 600
$ws.Range("F46").Value = 295
$ws.Range("F47").Value = 112
$ws.Range("F50").Value = 86
